# Add a "Date" column (H) to the import test sheet.
# H1 gets a real date value formatted as a short date; H2 is a blank date-
# formatted cell; H3, H6 and H7 hold bad (space) values but keep the date
# format; H4 and H5 hold the same bad value with no special format - this
# mirrors the "wrongnumeric" import-failure fixture used by the date-import
# tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: an actual date (2016-02-02), formatted as a short date (built-in
# numFmtId 14, "mm-dd-yy" maps onto it).
$dt = Get-Date -Year 2016 -Month 2 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("H1").Value = $dt.Date
$ws.Range("H1").NumberFormat = "mm-dd-yy"

# Reuse H1's freshly created style for the other date-formatted cells
# instead of re-applying NumberFormat (which would mint a duplicate style).
$ws.Range("H1").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H7").PasteSpecial(-4122)

# "Wrong" values - a lone space - simulating bad data for the date column.
$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

# Match the saved selection in the authored workbook.
$ws.Range("H7").Select() | Out-Null
